# SBA Calculator v1 example.xlsx
# Commit: "Change in foregiveness period"
#
# The "Covered Period" section on the Main sheet used to compute the
# forgiveness period automatically from two dates (Days = C8-C7, then
# Months = ROUNDDOWN((C8-C7)/30,1) -> 4.5). The author instead removed the
# "Days" helper row and hard-coded the forgiveness period (in months) to a
# literal value of 4, renaming the label to "Months - For Forgiveness".
# A threaded comment explaining the statutory 8-week covered period was
# also added on the (now repurposed) B10 cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

# 1. Remove the "Days" helper row (label in B9 and the =C8-C7 formula in C9).
$ws.Range("B9:C9").ClearContents()

# 2. Re-label B10 and hard-code C10 (the "Period" named range) to 4 months
#    instead of deriving it from ROUNDDOWN((C8-C7)/30,1).
$ws.Range("B10").Value = "Months - For Forgiveness"
$ws.Range("C10").Value = 4

# 3. Add a threaded comment on B10 documenting the statutory covered period.
$commentText = [char]0x28 + "3" + [char]0x29 + " the term " + [char]0x201C + "covered period" + [char]0x201D + " means the 8-week period beginning on the date of the origination of a covered loan;"
$ws.Range("B10").AddCommentThreaded($commentText) | Out-Null
